# se modif data para smoke en PreProd
$wb = $excel.ActiveWorkbook

# --- DatosCuenta ---
$wsCuenta = $wb.Worksheets.Item("DatosCuenta")
$wsCuenta.Range("A2").Value = "SmokePreTRES"
$wsCuenta.Range("B2").Value = "SmokePreNameTRES"
$wsCuenta.Range("D2").Value = 111
[void]$wsCuenta.Range("D2").Select()

# --- DatosHogar ---
$wsHogar = $wb.Worksheets.Item("DatosHogar")
$wsHogar.Range("A2").Value = 631

# --- DatosMotor ---
$wsMotor = $wb.Worksheets.Item("DatosMotor")
$wsMotor.Range("A2").Value = "SMA012"
$wsMotor.Range("B2").Value = "ABC12SSMA012"
$wsMotor.Range("C2").Value = "ZAZ123SSMA012"
[void]$wsMotor.Range("A4:C4").Select()

# --- DatosAP (stays the active sheet/tab) ---
$wsAP = $wb.Worksheets.Item("DatosAP")
$wsAP.Range("A2").Value = 21200111
[void]$wsAP.Range("A3").Select()
